$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.996.46"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").Value = "2.290.64"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.63"
$ws.Range("E5").Value = "  +1.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.47"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.530"
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("D9").Value = "2.309.71"
$ws.Range("E9").Value = "  +1.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0981"
$ws.Range("E10").Value = "  +3.00%  "
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.11"
$ws.Range("E12").Value = "  +8.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.341"
$ws.Range("E13").Value = "  +1.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.73"
$ws.Range("E14").Value = "  +4.29%  "
$ws.Range("D15").Value = "2.694.68"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("D16").Value = "55.019.44"
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000132"
$ws.Range("E17").Value = "  +1.90%  "
$ws.Range("D18").Value = "2.294.35"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.49"
$ws.Range("E19").Value = "  +2.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.17"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "312.73"
$ws.Range("E21").Value = "  +3.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.64"
$ws.Range("E22").Value = "  +4.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.41"
$ws.Range("E24").Value = "  -1.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.993"
$ws.Range("E25").Value = "  -0.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.153"
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.51"
$ws.Range("E27").Value = "  +3.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "172.24"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.16"
$ws.Range("E29").Value = "  +4.02%  "
$ws.Range("D30").Value = "0.0₃0709"
$ws.Range("E30").Value = "  +3.74%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.15"
$ws.Range("E31").Value = "  +6.81%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.64"
$ws.Range("E32").Value = "  +1.76%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.08"
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.993"
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.922"
$ws.Range("E36").Value = "  -3.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.23"
$ws.Range("E37").Value = "  +3.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.90"
$ws.Range("E38").Value = "  +5.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.85"
$ws.Range("E39").Value = "  +2.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.45"
$ws.Range("E40").Value = "  +3.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.376"
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "135.22"
$ws.Range("E42").Value = "  +8.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.11"
$ws.Range("E43").Value = "  +6.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.43"
$ws.Range("E44").Value = "  +2.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "261.00"
$ws.Range("E45").Value = "  +9.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0507"
$ws.Range("E46").Value = "  +2.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0912"
$ws.Range("E47").Value = "  +2.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.552"
$ws.Range("E48").Value = "  +1.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.376"
$ws.Range("E49").Value = "  +1.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0212"
$ws.Range("E50").Value = "  +3.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.48"
$ws.Range("E51").Value = "  +1.70%  "
